$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 1693
$ws1.Range("F11").Value = 1686
$ws1.Range("F13").Value = 87
$ws1.Range("F16").Value = 194
$ws1.Range("F21").Value = 275
$ws1.Range("F24").Value = 229

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1693
$ws4.Range("F12").Value = 1686
$ws4.Range("F14").Value = 87
$ws4.Range("F17").Value = 194
$ws4.Range("F22").Value = 275
$ws4.Range("F25").Value = 229
